# "Nov 8th - Status" update to the Daily_Status workbook.
# Appends three more days of status entries (06/11, 07/11, 08/11) under the
# existing "DATE / Work Status / ..." log on Sheet1, right after row 84.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- 06/11/2021 : holiday -------------------------------------------------
$ws.Range("A85").Value = "06/11/2021"
$ws.Range("B85").Value = "HOLIDAY"

# --- 07/11/2021 : holiday -------------------------------------------------
$ws.Range("A86").Value = "07/11/2021"
$ws.Range("B86").Value = "HOLIDAY"

# --- 08/11/2021 : status notes --------------------------------------------
$ws.Range("A87").Value = "08/11/2021"
$ws.Range("B87").Value = "Continued on codec 2.0 "
$ws.Range("C87").Value = "Updating the notes"

$ws.Range("B88").Value = "Reported the current study progress with refrerence and links"
$ws.Range("B89").Value = "Attended the session on LDD recap"
$ws.Range("B90").Value = "Had internal discussion with the teammates"

# Row 77 wraps onto one fewer line in the refreshed layout.
$ws.Rows.Item(77).RowHeight = 30

# Scroll the window down to the newly-added rows and leave the cursor on
# the last entry, same as the author did before saving.
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B90").Select() | Out-Null
